# Auto-generated Excel COM-interop script to refresh market-data derived
# columns (H..N) across multiple worksheets, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3667.125
$ws.Range("I70").Value = 1734.5
$ws.Range("J70").Value = 4311.3335
$ws.Range("K70").Value = 5203.5
$ws.Range("L70").Value = 12934.0005
$ws.Range("M70").Value = -4933.5
$ws.Range("N70").Value = -13474.0005
$ws.Range("H73").Value = 3667.125
$ws.Range("I73").Value = 1734.5
$ws.Range("J73").Value = 4311.3335
$ws.Range("K73").Value = 5203.5
$ws.Range("L73").Value = 12934.0005
$ws.Range("M73").Value = -4267.5
$ws.Range("N73").Value = -14806.0005
$ws.Range("H74").Value = 3779.6
$ws.Range("I74").Value = 3299.6667
$ws.Range("J74").Value = 4499.5
$ws.Range("K74").Value = 3299.6667
$ws.Range("L74").Value = 4499.5
$ws.Range("M74").Value = -2363.6667
$ws.Range("N74").Value = -6371.5
$ws.Range("H77").Value = 3779.6
$ws.Range("I77").Value = 3299.6667
$ws.Range("J77").Value = 4499.5
$ws.Range("K77").Value = 16498.3335
$ws.Range("L77").Value = 22497.5
$ws.Range("M77").Value = -11818.3335
$ws.Range("N77").Value = -31857.5
$ws.Range("H107").Value = 1367.0476
$ws.Range("I107").Value = 1000.82355
$ws.Range("K107").Value = 1000.82355
$ws.Range("M107").Value = 919.17645
$ws.Range("H132").Value = 2138.276
$ws.Range("I132").Value = 2281.7727
$ws.Range("J132").Value = 1687.2858
$ws.Range("K132").Value = 6845.3181
$ws.Range("L132").Value = 5061.857400000001
$ws.Range("M132").Value = -4315.3181
$ws.Range("N132").Value = -10121.8574
$ws.Range("H137").Value = 1481.76
$ws.Range("I137").Value = 1097.5
$ws.Range("K137").Value = 3292.5
$ws.Range("M137").Value = -742.5
$ws.Range("H138").Value = 6096.926
$ws.Range("J138").Value = 8170.4375
$ws.Range("L138").Value = 24511.3125
$ws.Range("N138").Value = -34791.3125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1804.1086
$ws.Range("I32").Value = 1053.8942
$ws.Range("K32").Value = 1053.8942
$ws.Range("M32").Value = -766.8942
$ws.Range("H61").Value = 31254972
$ws.Range("I61").Value = 34487704
$ws.Range("J61").Value = 5232.3335
$ws.Range("K61").Value = 34487704
$ws.Range("L61").Value = 5232.3335
$ws.Range("M61").Value = -34487492
$ws.Range("N61").Value = -5656.3335
$ws.Range("H110").Value = 78847.92
$ws.Range("I110").Value = 112112.336
$ws.Range("K110").Value = 112112.336
$ws.Range("M110").Value = -110067.336
$ws.Range("H122").Value = 6730.533
$ws.Range("I122").Value = 5573.6924
$ws.Range("K122").Value = 16721.0772
$ws.Range("M122").Value = -14271.0772
$ws.Range("H132").Value = 2943614.5
$ws.Range("I132").Value = 3228112.8
$ws.Range("J132").Value = 3799
$ws.Range("K132").Value = 9684338.399999999
$ws.Range("L132").Value = 11397
$ws.Range("M132").Value = -9681808.399999999
$ws.Range("N132").Value = -16457
$ws.Range("H136").Value = 31254972
$ws.Range("I136").Value = 34487704
$ws.Range("J136").Value = 5232.3335
$ws.Range("K136").Value = 103463112
$ws.Range("L136").Value = 15697.0005
$ws.Range("M136").Value = -103460562
$ws.Range("N136").Value = -20797.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 126076.81
$ws.Range("I107").Value = 524.2222
$ws.Range("K107").Value = 524.2222
$ws.Range("M107").Value = 1395.7778
$ws.Range("H134").Value = 38464220
$ws.Range("I134").Value = 41669416
$ws.Range("J134").Value = 1899
$ws.Range("K134").Value = 125008248
$ws.Range("L134").Value = 5697
$ws.Range("M134").Value = -125005713
$ws.Range("N134").Value = -10767
$ws.Range("H135").Value = 89999
$ws.Range("J135").Value = 89999
$ws.Range("L135").Value = 89999
$ws.Range("N135").Value = -100139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 945
$ws.Range("I16").Value = 940
$ws.Range("J16").Value = 950
$ws.Range("K16").Value = 940
$ws.Range("L16").Value = 950
$ws.Range("M16").Value = -653
$ws.Range("N16").Value = -1524
$ws.Range("H22").Value = 547
$ws.Range("I22").Value = 558.2273
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 558.2273
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -208.2273
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 7625.8335
$ws.Range("I31").Value = 11009.857
$ws.Range("K31").Value = 11009.857
$ws.Range("M31").Value = -10714.857
$ws.Range("H34").Value = 7625.8335
$ws.Range("I34").Value = 11009.857
$ws.Range("K34").Value = 11009.857
$ws.Range("M34").Value = -10807.857
$ws.Range("H58").Value = 17248316
$ws.Range("I58").Value = 22735212
$ws.Range("J58").Value = 3783.2856
$ws.Range("K58").Value = 22735212
$ws.Range("L58").Value = 3783.2856
$ws.Range("M58").Value = -22735009
$ws.Range("N58").Value = -4189.2856
$ws.Range("H113").Value = 945
$ws.Range("I113").Value = 940
$ws.Range("J113").Value = 950
$ws.Range("K113").Value = 940
$ws.Range("L113").Value = 950
$ws.Range("M113").Value = 1230
$ws.Range("N113").Value = -5290
$ws.Range("H132").Value = 52634904
$ws.Range("I132").Value = 66669744
$ws.Range("K132").Value = 200009232
$ws.Range("M132").Value = -200006702
$ws.Range("H134").Value = 8066522.5
$ws.Range("I134").Value = 8930510
$ws.Range("K134").Value = 26791530
$ws.Range("M134").Value = -26788995
$ws.Range("H136").Value = 17248316
$ws.Range("I136").Value = 22735212
$ws.Range("J136").Value = 3783.2856
$ws.Range("K136").Value = 68205636
$ws.Range("L136").Value = 11349.8568
$ws.Range("M136").Value = -68203086
$ws.Range("N136").Value = -16449.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4914119
$ws.Range("I4").Value = 5381654
$ws.Range("K4").Value = 16144962
$ws.Range("M4").Value = -16144850
$ws.Range("H12").Value = 115.5
$ws.Range("I12").Value = 15.166667
$ws.Range("K12").Value = 45.500001
$ws.Range("M12").Value = 127.499999
$ws.Range("H98").Value = 2316.5
$ws.Range("J98").Value = 2207.3333
$ws.Range("L98").Value = 6621.999899999999
$ws.Range("N98").Value = -9617.999899999999
$ws.Range("H137").Value = 1831.75
$ws.Range("I137").Value = 1593.4286
$ws.Range("K137").Value = 4780.2858
$ws.Range("M137").Value = 319.7142000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 30000
$ws.Range("J51").Value = 30000
$ws.Range("L51").Value = 30000
$ws.Range("N51").Value = -31018
$ws.Range("H80").Value = 1932.3334
$ws.Range("I80").Value = 1948.5
$ws.Range("J80").Value = 1900
$ws.Range("K80").Value = 1948.5
$ws.Range("L80").Value = 1900
$ws.Range("M80").Value = -950.5
$ws.Range("N80").Value = -3896
$ws.Range("H83").Value = 1932.3334
$ws.Range("I83").Value = 1948.5
$ws.Range("J83").Value = 1900
$ws.Range("K83").Value = 9742.5
$ws.Range("L83").Value = 9500
$ws.Range("M83").Value = -4750.5
$ws.Range("N83").Value = -19484
$ws.Range("H94").Value = 29499.5
$ws.Range("J94").Value = 29499.5
$ws.Range("L94").Value = 29499.5
$ws.Range("N94").Value = -30851.5
$ws.Range("H113").Value = 51252.715
$ws.Range("I113").Value = 59439.5
$ws.Range("J113").Value = 2132
$ws.Range("K113").Value = 59439.5
$ws.Range("L113").Value = 2132
$ws.Range("M113").Value = -57269.5
$ws.Range("N113").Value = -6472
$ws.Range("H122").Value = 3101.423
$ws.Range("I122").Value = 1571.3572
$ws.Range("J122").Value = 4886.5
$ws.Range("K122").Value = 4714.071599999999
$ws.Range("L122").Value = 14659.5
$ws.Range("M122").Value = -2264.071599999999
$ws.Range("N122").Value = -19559.5
$ws.Range("H126").Value = 2199
$ws.Range("I126").Value = 1699
$ws.Range("J126").Value = 2699
$ws.Range("K126").Value = 5097
$ws.Range("L126").Value = 8097
$ws.Range("M126").Value = -2627
$ws.Range("N126").Value = -13037

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20874394
$ws.Range("I132").Value = 22862260
$ws.Range("J132").Value = 1789.5
$ws.Range("K132").Value = 68586780
$ws.Range("L132").Value = 5368.5
$ws.Range("M132").Value = -68584250
$ws.Range("N132").Value = -10428.5
$ws.Range("H136").Value = 6185
$ws.Range("I136").Value = 4872
$ws.Range("K136").Value = 14616
$ws.Range("M136").Value = -12066

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 23810740
$ws.Range("I132").Value = 26316972
$ws.Range("K132").Value = 78950916
$ws.Range("M132").Value = -78948386
$ws.Range("H136").Value = 55556320
$ws.Range("I136").Value = 55556320
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 166668960
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -166666410
$ws.Range("N136").ClearContents()

